$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 442.36365
$ws.Range("I33").Value = 214.11765
$ws.Range("J33").Value = 1218.4
$ws.Range("K33").Value = 214.11765
$ws.Range("L33").Value = 1218.4
$ws.Range("M33").Value = 14.88235
$ws.Range("N33").Value = -1676.4
$ws.Range("H113").Value = 2058.9546
$ws.Range("I113").Value = 2019.8
$ws.Range("J113").Value = 2142.8572
$ws.Range("K113").Value = 2019.8
$ws.Range("L113").Value = 2142.8572
$ws.Range("M113").Value = 1234.2
$ws.Range("N113").Value = -8650.8572
$ws.Range("H132").Value = 191482.27
$ws.Range("I132").Value = 2730.175
$ws.Range("K132").Value = 8190.525000000001
$ws.Range("M132").Value = -5660.525000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5683.35
$ws.Range("I32").Value = 5369.3906
$ws.Range("J32").Value = 7784.4614
$ws.Range("K32").Value = 5369.3906
$ws.Range("L32").Value = 7784.4614
$ws.Range("M32").Value = -5082.3906
$ws.Range("N32").Value = -8358.4614
$ws.Range("H61").Value = 2075.4187
$ws.Range("I61").Value = 1969.5526
$ws.Range("J61").Value = 2880
$ws.Range("K61").Value = 1969.5526
$ws.Range("L61").Value = 2880
$ws.Range("M61").Value = -1757.5526
$ws.Range("N61").Value = -3304
$ws.Range("H136").Value = 2075.4187
$ws.Range("I136").Value = 1969.5526
$ws.Range("J136").Value = 2880
$ws.Range("K136").Value = 5908.6578
$ws.Range("L136").Value = 8640
$ws.Range("M136").Value = -3358.6578
$ws.Range("N136").Value = -13740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 43524410
$ws.Range("I134").Value = 2754
$ws.Range("J134").Value = 125127500
$ws.Range("K134").Value = 8262
$ws.Range("L134").Value = 375382500
$ws.Range("M134").Value = -5727
$ws.Range("N134").Value = -375387570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5144.278
$ws.Range("I31").Value = 1635.2858
$ws.Range("J31").Value = 7377.273
$ws.Range("K31").Value = 1635.2858
$ws.Range("L31").Value = 7377.273
$ws.Range("M31").Value = -1340.2858
$ws.Range("N31").Value = -7967.273
$ws.Range("H34").Value = 5144.278
$ws.Range("I34").Value = 1635.2858
$ws.Range("J34").Value = 7377.273
$ws.Range("K34").Value = 1635.2858
$ws.Range("L34").Value = 7377.273
$ws.Range("M34").Value = -1433.2858
$ws.Range("N34").Value = -7781.273
$ws.Range("H58").Value = 1580.762
$ws.Range("I58").Value = 1142.9375
$ws.Range("J58").Value = 2981.8
$ws.Range("K58").Value = 1142.9375
$ws.Range("L58").Value = 2981.8
$ws.Range("M58").Value = -939.9375
$ws.Range("N58").Value = -3387.8
$ws.Range("H99").Value = 3783.7036
$ws.Range("I99").Value = 3161.875
$ws.Range("J99").Value = 4688.1816
$ws.Range("K99").Value = 3161.875
$ws.Range("L99").Value = 4688.1816
$ws.Range("M99").Value = -1663.875
$ws.Range("N99").Value = -7684.1816
$ws.Range("H126").Value = 3783.7036
$ws.Range("I126").Value = 3161.875
$ws.Range("J126").Value = 4688.1816
$ws.Range("K126").Value = 9485.625
$ws.Range("L126").Value = 14064.5448
$ws.Range("M126").Value = -7015.625
$ws.Range("N126").Value = -19004.5448
$ws.Range("H132").Value = 114145.336
$ws.Range("I132").Value = 251703
$ws.Range("J132").Value = 4099.2
$ws.Range("K132").Value = 755109
$ws.Range("L132").Value = 12297.6
$ws.Range("M132").Value = -752579
$ws.Range("N132").Value = -17357.6
$ws.Range("H134").Value = 1640.9688
$ws.Range("I134").Value = 1482.5358
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 4447.607400000001
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -1912.607400000001
$ws.Range("N134").Value = -13320
$ws.Range("H136").Value = 1580.762
$ws.Range("I136").Value = 1142.9375
$ws.Range("J136").Value = 2981.8
$ws.Range("K136").Value = 3428.8125
$ws.Range("L136").Value = 8945.400000000001
$ws.Range("M136").Value = -878.8125
$ws.Range("N136").Value = -14045.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 22591.4
$ws.Range("J100").Value = 27989.25
$ws.Range("L100").Value = 83967.75
$ws.Range("N100").Value = -85589.75
$ws.Range("H113").Value = 495
$ws.Range("I113").Value = 497.11765
$ws.Range("J113").Value = 486
$ws.Range("K113").Value = 1491.35295
$ws.Range("L113").Value = 1458
$ws.Range("M113").Value = 678.64705
$ws.Range("N113").Value = -5798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1872.2963
$ws.Range("I122").Value = 1636.625
$ws.Range("J122").Value = 2215.0908
$ws.Range("K122").Value = 4909.875
$ws.Range("L122").Value = 6645.2724
$ws.Range("M122").Value = -2459.875
$ws.Range("N122").Value = -11545.2724
$ws.Range("H126").Value = 1722.2572
$ws.Range("I126").Value = 1436.9231
$ws.Range("J126").Value = 1890.8636
$ws.Range("K126").Value = 4310.7693
$ws.Range("L126").Value = 5672.5908
$ws.Range("M126").Value = -1840.7693
$ws.Range("N126").Value = -10612.5908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2514.8286
$ws.Range("I7").Value = 2527.8845
$ws.Range("J7").Value = 2477.111
$ws.Range("K7").Value = 2527.8845
$ws.Range("L7").Value = 2477.111
$ws.Range("M7").Value = -2415.8845
$ws.Range("N7").Value = -2701.111
$ws.Range("H40").Value = 48476
$ws.Range("I40").Value = 2135.3333
$ws.Range("J40").Value = 65853.75
$ws.Range("K40").Value = 2135.3333
$ws.Range("L40").Value = 65853.75
$ws.Range("M40").Value = -1999.3333
$ws.Range("N40").Value = -66125.75
$ws.Range("H126").Value = 2514.8286
$ws.Range("I126").Value = 2527.8845
$ws.Range("J126").Value = 2477.111
$ws.Range("K126").Value = 7583.6535
$ws.Range("L126").Value = 7431.333
$ws.Range("M126").Value = -5113.6535
$ws.Range("N126").Value = -12371.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1140.6
$ws.Range("I126").Value = 1266
$ws.Range("K126").Value = 3798
$ws.Range("M126").Value = -1328
$ws.Range("H136").Value = 1567595.5
$ws.Range("I136").Value = 1881176.2
$ws.Range("J136").Value = 716447.9
$ws.Range("K136").Value = 5643528.6
$ws.Range("L136").Value = 2149343.7
$ws.Range("M136").Value = -5640978.6
$ws.Range("N136").Value = -2154443.7
